$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.037.56"
$ws.Range("E2").Value = "  +4.97%  "
$ws.Range("D3").Value = "2.779.86"
$ws.Range("E3").Value = "  +5.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'340.40"
$ws.Range("E5").Value = "  +4.22%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'115.36"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("E7").Value = "  +4.85%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.576"
$ws.Range("E9").Value = "  +4.80%  "
$ws.Range("D10").Value = "'41.81"
$ws.Range("E10").Value = "  +5.42%  "
$ws.Range("E11").Value = "  +5.27%  "
$ws.Range("D12").Value = "'20.05"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "'7.60"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "3.214.39"
$ws.Range("E15").Value = "  +5.00%  "
$ws.Range("D16").Value = "2.773.65"
$ws.Range("E16").Value = "  +5.08%  "
$ws.Range("D17").Value = "51.868.18"
$ws.Range("E17").Value = "  +4.63%  "
$ws.Range("D18").Value = "'0.877"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("E19").Value = "  +9.92%  "
$ws.Range("D20").Value = "'6.99"
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("D21").Value = "'13.22"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D23").Value = "'276.15"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "'69.89"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  +7.11%  "
$ws.Range("D26").Value = "'26.69"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'34.68"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'50.09"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'5.71"
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("D34").Value = "'0.0819"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'18.98"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'2.10"
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("D38").Value = "'4.94"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "'3.21"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").Value = "'0.0379"
$ws.Range("E40").Value = "  +8.92%  "
$ws.Range("E41").Value = "  +27.44%  "
$ws.Range("D42").Value = "'2.36"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("E43").Value = "  +3.08%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'23.19"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'125.71"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "2.067.97"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "'5.55"
$ws.Range("E49").Value = "  +5.78%  "
$ws.Range("D50").Value = "'8.84"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "'0.881"
$ws.Range("E51").Value = "  +14.11%  "

# Reset style on cells that were forced to text via a leading apostrophe,
# so they don't retain an Excel 'Number Stored as Text' quote-prefix style.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
